# Update currency data sheet:
#  - rename header cells E1/F1, add new header G1
#  - refresh timestamp, change-% (now text), trend label, and add icon column
#    for the first 4 data rows
#  - remove the now-duplicate rows 6-9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("E1").Value = "Change_Pct"
$ws.Range("F1").Value = "Trend"
$ws.Range("G1").Value = "Icon"
# Give the new header cell the same look (bold/border/center) as the rest
# of row 1 by copying F1's formatting over.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

# --- Remove the stale duplicate rows (old rows 6-9) first, so the ----
# --- remaining data rows keep indices 2-5 -----------------------------
$ws.Range("A6:G9").EntireRow.Delete()

# --- Refresh rows 2-5 ---------------------------------------------------
$newTimestamp = "08:09:00"
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
    # Change_% column becomes a text value instead of a number. Prefixing
    # with an apostrophe forces text storage; resetting the style back to
    # Normal drops the quote-prefix formatting Excel would otherwise apply.
    $changeText = $ws.Cells.Item($r, 5).Text
    $changeCell = $ws.Cells.Item($r, 5)
    $changeCell.Value = "'" + $changeText
    $changeCell.Style = "Normal"
    $ws.Cells.Item($r, 6).Value = "ALTA"
    $ws.Cells.Item($r, 7).Value = "🟢"
}
